$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.906.71'
$ws.Range("E2").Value = '  +0.08%  '

# Row 3
$ws.Range("D3").Value = '1.895.15'
$ws.Range("E3").Value = '  -0.08%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7730'
$ws.Range("E5").Value = '  -2.54%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.74'
$ws.Range("E6").Value = '  +0.42%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3133'
$ws.Range("E8").Value = '  -0.83%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.69'
$ws.Range("E9").Value = '  +1.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07261'
$ws.Range("E10").Value = '  +0.33%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08894'
$ws.Range("E11").Value = '  +9.70%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7725'
$ws.Range("E12").Value = '  +0.75%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.444'
$ws.Range("E13").Value = '  -2.77%  '

# Row 14
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.40'
$ws.Range("E14").Value = '  +1.95%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.827.55'
$ws.Range("E15").Value = '  -3.23%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.199'
$ws.Range("E16").Value = '  +0.12%  '

# Row 17
$ws.Range("D17").Value = '29.843.19'
$ws.Range("E17").Value = '  -0.14%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.97'
$ws.Range("E18").Value = '  +0.09%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.92'
$ws.Range("E19").Value = '  +0.62%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007865'
$ws.Range("E20").Value = '  +0.78%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.135'
$ws.Range("E21").Value = '  -1.04%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.03%  '

# Row 23
$ws.Range("D23").Value = '2.093.69'
$ws.Range("E23").Value = '  -0.76%  '

# Row 24
$ws.Range("E24").Value = '  +0.03%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1600'
$ws.Range("E25").Value = '  -4.34%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.513'
$ws.Range("E26").Value = '  +0.78%  '

# Row 27
$ws.Range("E27").Value = '  -0.76%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.83'
$ws.Range("E28").Value = '  +0.61%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.046'
$ws.Range("E29").Value = '  -1.10%  '

# Row 30
$ws.Range("E30").Value = '  +1.86%  '

# Row 31
$ws.Range("E31").Value = '  -0.34%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.546'
$ws.Range("E32").Value = '  +1.49%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.105'
$ws.Range("E33").Value = '  +0.27%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05512'
$ws.Range("E34").Value = '  -0.13%  '

# Row 35
$ws.Range("E35").Value = '  -2.51%  '

# Row 36
$ws.Range("E36").Value = '  +1.53%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9947'
$ws.Range("E37").Value = '  +0.02%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.719'
$ws.Range("E38").Value = '  +3.77%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01965'
$ws.Range("E39").Value = '  +1.97%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.790'
$ws.Range("E40").Value = '  +0.24%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4503'
$ws.Range("E41").Value = '  +1.85%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '73.68'
$ws.Range("E42").Value = '  -1.01%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.048'
$ws.Range("E43").Value = '  +2.65%  '

# Row 44
$ws.Range("D44").Value = '1.087.77'
$ws.Range("E44").Value = '  -6.05%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8556'
$ws.Range("E45").Value = '  +0.39%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  +0.06%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.889'
$ws.Range("E47").Value = '  +0.60%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.62'
$ws.Range("E48").Value = '  -2.08%  '

# Row 49
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.596'
$ws.Range("E49").Value = '  +1.83%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.851'
$ws.Range("E50").Value = '  -1.71%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.991'
$ws.Range("E51").Value = '  -1.91%  '
